$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 161
$ws.Cells.Item(161, 2).Value = 57756
$ws.Cells.Item(161, 4).Value = 66.44
$ws.Cells.Item(161, 5).Value = 79.37
$ws.Cells.Item(161, 6).Value = -100
$ws.Cells.Item(161, 7).Value = -6644

# Row 163
$ws.Cells.Item(163, 2).Value = 53925
$ws.Cells.Item(163, 4).Value = 66.44
$ws.Cells.Item(163, 5).Value = 79.37
$ws.Cells.Item(163, 6).Value = 1
$ws.Cells.Item(163, 7).Value = 66.44

# Row 183
$ws.Cells.Item(183, 2).Value = 57552
$ws.Cells.Item(183, 4).Value = 120.69
$ws.Cells.Item(183, 5).Value = 136.86
$ws.Cells.Item(183, 6).Value = -5
$ws.Cells.Item(183, 7).Value = -603.45

# Row 184
$ws.Cells.Item(184, 2).Value = 64329
$ws.Cells.Item(184, 4).Value = 120.69
$ws.Cells.Item(184, 5).Value = 128.32
$ws.Cells.Item(184, 6).Value = 6
$ws.Cells.Item(184, 7).Value = 724.14

# Row 264
$ws.Cells.Item(264, 2).Value = 48719
$ws.Cells.Item(264, 4).Value = 295.75
$ws.Cells.Item(264, 5).Value = 353.35
$ws.Cells.Item(264, 6).Value = -81
$ws.Cells.Item(264, 7).Value = -23955.75

# Row 265
$ws.Cells.Item(265, 2).Value = 64979
$ws.Cells.Item(265, 4).Value = 295.75
$ws.Cells.Item(265, 5).Value = 314.41
$ws.Cells.Item(265, 6).Value = 82
$ws.Cells.Item(265, 7).Value = 24251.5

# Row 279
$ws.Cells.Item(279, 2).Value = 48706
$ws.Cells.Item(279, 4).Value = 33.3
$ws.Cells.Item(279, 5).Value = 39.8
$ws.Cells.Item(279, 6).Value = -144
$ws.Cells.Item(279, 7).Value = -4795.2

# Row 280
$ws.Cells.Item(280, 2).Value = 64973
$ws.Cells.Item(280, 4).Value = 33.3
$ws.Cells.Item(280, 5).Value = 35.4
$ws.Cells.Item(280, 6).Value = 150
$ws.Cells.Item(280, 7).Value = 4995

# Row 316
$ws.Cells.Item(316, 2).Value = 61610
$ws.Cells.Item(316, 4).Value = 102.71
$ws.Cells.Item(316, 5).Value = 122.71
$ws.Cells.Item(316, 6).Value = -58
$ws.Cells.Item(316, 7).Value = -5957.18

# Row 317
$ws.Cells.Item(317, 2).Value = 63565
$ws.Cells.Item(317, 4).Value = 102.71
$ws.Cells.Item(317, 5).Value = 109.19
$ws.Cells.Item(317, 6).Value = 60
$ws.Cells.Item(317, 7).Value = 6162.6

# Row 318
$ws.Cells.Item(318, 2).Value = 57077
$ws.Cells.Item(318, 4).Value = 93.08
$ws.Cells.Item(318, 5).Value = 111.2
$ws.Cells.Item(318, 6).Value = 1
$ws.Cells.Item(318, 7).Value = 93.08

# Row 346
$ws.Cells.Item(346, 2).Value = 63520
$ws.Cells.Item(346, 4).Value = 144.28
$ws.Cells.Item(346, 5).Value = 153.4
$ws.Cells.Item(346, 6).Value = 97
$ws.Cells.Item(346, 7).Value = 13995.16

# Row 347
$ws.Cells.Item(347, 2).Value = 55373
$ws.Cells.Item(347, 4).Value = 144.28
$ws.Cells.Item(347, 5).Value = 163.62
$ws.Cells.Item(347, 6).Value = -94
$ws.Cells.Item(347, 7).Value = -13562.32

# Row 350
$ws.Cells.Item(350, 2).Value = 63571
$ws.Cells.Item(350, 4).Value = 143.48
$ws.Cells.Item(350, 5).Value = 152.53
$ws.Cells.Item(350, 6).Value = 29
$ws.Cells.Item(350, 7).Value = 4160.92

# Row 351
$ws.Cells.Item(351, 2).Value = 63531
$ws.Cells.Item(351, 4).Value = 143.48
$ws.Cells.Item(351, 5).Value = 152.53
$ws.Cells.Item(351, 6).Value = 80
$ws.Cells.Item(351, 7).Value = 11478.4

# Row 352
$ws.Cells.Item(352, 2).Value = 57802
$ws.Cells.Item(352, 4).Value = 143.48
$ws.Cells.Item(352, 5).Value = 162.71
$ws.Cells.Item(352, 6).Value = -79
$ws.Cells.Item(352, 7).Value = -11334.92

# Row 379
$ws.Cells.Item(379, 2).Value = 61608
$ws.Cells.Item(379, 4).Value = 129.01
$ws.Cells.Item(379, 5).Value = 154.12
$ws.Cells.Item(379, 6).Value = -56
$ws.Cells.Item(379, 7).Value = -7224.56

# Row 380
$ws.Cells.Item(380, 2).Value = 63564
$ws.Cells.Item(380, 4).Value = 129.01
$ws.Cells.Item(380, 5).Value = 137.16
$ws.Cells.Item(380, 6).Value = 57
$ws.Cells.Item(380, 7).Value = 7353.57

# Row 382
$ws.Cells.Item(382, 2).Value = 63560
$ws.Cells.Item(382, 4).Value = 126.86
$ws.Cells.Item(382, 5).Value = 134.87
$ws.Cells.Item(382, 6).Value = 104
$ws.Cells.Item(382, 7).Value = 13193.44

# Row 383
$ws.Cells.Item(383, 2).Value = 60325
$ws.Cells.Item(383, 4).Value = 126.86
$ws.Cells.Item(383, 5).Value = 151.57
$ws.Cells.Item(383, 6).Value = -102
$ws.Cells.Item(383, 7).Value = -12939.72

# Row 457
$ws.Cells.Item(457, 2).Value = 63681
$ws.Cells.Item(457, 4).Value = 22.42
$ws.Cells.Item(457, 5).Value = 23.84
$ws.Cells.Item(457, 6).Value = 65
$ws.Cells.Item(457, 7).Value = 1457.3

# Row 458
$ws.Cells.Item(458, 2).Value = 31930
$ws.Cells.Item(458, 4).Value = 22.42
$ws.Cells.Item(458, 5).Value = 26.8
$ws.Cells.Item(458, 6).Value = -62
$ws.Cells.Item(458, 7).Value = -1390.04

# Row 536
$ws.Cells.Item(536, 2).Value = 58047
$ws.Cells.Item(536, 4).Value = 105.54
$ws.Cells.Item(536, 5).Value = 126.1
$ws.Cells.Item(536, 6).Value = 54
$ws.Cells.Item(536, 7).Value = 5699.16

# Row 537
$ws.Cells.Item(537, 2).Value = 47097
$ws.Cells.Item(537, 4).Value = 112.28
$ws.Cells.Item(537, 5).Value = 134.16
$ws.Cells.Item(537, 6).Value = 15
$ws.Cells.Item(537, 7).Value = 1684.2

# Row 581
$ws.Cells.Item(581, 2).Value = 65068
$ws.Cells.Item(581, 4).Value = 13.15
$ws.Cells.Item(581, 5).Value = 13.97
$ws.Cells.Item(581, 6).Value = 232
$ws.Cells.Item(581, 7).Value = 3050.8

# Row 582
$ws.Cells.Item(582, 2).Value = 53602
$ws.Cells.Item(582, 4).Value = 13.15
$ws.Cells.Item(582, 5).Value = 15.69
$ws.Cells.Item(582, 6).Value = -231
$ws.Cells.Item(582, 7).Value = -3037.65

# Row 586
$ws.Cells.Item(586, 2).Value = 45695
$ws.Cells.Item(586, 4).Value = 19.73
$ws.Cells.Item(586, 5).Value = 23.58
$ws.Cells.Item(586, 6).Value = -36
$ws.Cells.Item(586, 7).Value = -710.28

# Row 587
$ws.Cells.Item(587, 2).Value = 64915
$ws.Cells.Item(587, 4).Value = 19.73
$ws.Cells.Item(587, 5).Value = 20.98
$ws.Cells.Item(587, 6).Value = 40
$ws.Cells.Item(587, 7).Value = 789.2

# Row 590
$ws.Cells.Item(590, 2).Value = 45706
$ws.Cells.Item(590, 4).Value = 19.73
$ws.Cells.Item(590, 5).Value = 23.58
$ws.Cells.Item(590, 6).Value = -202
$ws.Cells.Item(590, 7).Value = -3985.46

# Row 591
$ws.Cells.Item(591, 2).Value = 64922
$ws.Cells.Item(591, 4).Value = 19.73
$ws.Cells.Item(591, 5).Value = 20.98
$ws.Cells.Item(591, 6).Value = 207
$ws.Cells.Item(591, 7).Value = 4084.11

# Row 599
$ws.Cells.Item(599, 2).Value = 45709
$ws.Cells.Item(599, 4).Value = 13.15
$ws.Cells.Item(599, 5).Value = 15.69
$ws.Cells.Item(599, 6).Value = -300
$ws.Cells.Item(599, 7).Value = -3945

# Row 600
$ws.Cells.Item(600, 2).Value = 64925
$ws.Cells.Item(600, 4).Value = 13.15
$ws.Cells.Item(600, 5).Value = 13.97
$ws.Cells.Item(600, 6).Value = 302
$ws.Cells.Item(600, 7).Value = 3971.3

# Row 601
$ws.Cells.Item(601, 2).Value = 45702
$ws.Cells.Item(601, 4).Value = 26.3
$ws.Cells.Item(601, 5).Value = 31.43
$ws.Cells.Item(601, 6).Value = -215
$ws.Cells.Item(601, 7).Value = -5654.5

# Row 602
$ws.Cells.Item(602, 2).Value = 64919
$ws.Cells.Item(602, 4).Value = 26.3
$ws.Cells.Item(602, 5).Value = 27.97
$ws.Cells.Item(602, 6).Value = 224
$ws.Cells.Item(602, 7).Value = 5891.2

# Row 687
$ws.Cells.Item(687, 2).Value = 53319
$ws.Cells.Item(687, 4).Value = 273.92
$ws.Cells.Item(687, 5).Value = 310.64
$ws.Cells.Item(687, 6).Value = -6
$ws.Cells.Item(687, 7).Value = -1643.52

# Row 688
$ws.Cells.Item(688, 2).Value = 64810
$ws.Cells.Item(688, 4).Value = 273.92
$ws.Cells.Item(688, 5).Value = 291.22
$ws.Cells.Item(688, 6).Value = 7
$ws.Cells.Item(688, 7).Value = 1917.44

# Row 720
$ws.Cells.Item(720, 2).Value = 60022
$ws.Cells.Item(720, 4).Value = 32.83
$ws.Cells.Item(720, 5).Value = 37.22
$ws.Cells.Item(720, 6).Value = -113
$ws.Cells.Item(720, 7).Value = -3709.79

# Row 721
$ws.Cells.Item(721, 2).Value = 64830
$ws.Cells.Item(721, 4).Value = 32.83
$ws.Cells.Item(721, 5).Value = 34.9
$ws.Cells.Item(721, 6).Value = 117
$ws.Cells.Item(721, 7).Value = 3841.11
